$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34 status changes from "Reminder email" to "Sent key" (reuses an
# existing shared string).
$ws.Range("D34").Value = "Sent key"

# --- New rows 44 and 45 appended to the "name / key" list at the bottom.
$ws.Range("A44").Value = "TOT game"
$ws.Range("M44").Value = "KMDIG-99F0K-T0BLE"

$ws.Range("A45").Value = "Tom Johnson"
$ws.Range("M45").Value = " YMVKC-VGL02-RCXAL"

# --- Update the "State" column (D) wording: "Made positive video" / "Made
# --- negative video" became "Created video" / "Created negative video".
# --- Rows that previously said "Reminder email" but now indicate the video
# --- was created also get updated to "Created video".
$ws.Range("D2").Value = "Created video"
$ws.Range("D3").Value = "Created video"
$ws.Range("D4").Value = "Created video"
$ws.Range("D8").Value = "Created video"
$ws.Range("D9").Value = "Created video"
$ws.Range("D19").Value = "Created video"
$ws.Range("D20").Value = "Created negative video"

# Row 34 now also has a key listed.
$ws.Range("M34").Value = "0N2D5-62D32-WP7R5"

# Row 45 is slightly taller and M45 uses a distinct font (Arial 12,
# color #222222).
$ws.Range("A45:M45").RowHeight = 15.75
$ws.Range("M45").Font.Name = "Arial"
$ws.Range("M45").Font.Size = 12
$ws.Range("M45").Font.Color = 2236962

# Move the active selection to M34 (single cell) to match the saved view
# state.
$ws.Range("M34").Select()
